# Agrupamiento de datos por semana
# Consolidate the duplicated "date" rows (one row carrying the federal/state
# tax amounts, one row carrying the formatted 941/EDD payment amounts) into a
# single row per week, dropping the extra duplicate rows for 11152024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$companyName = "ALEJANDRA S FASHIÓON INC"

# New consolidated data, one entry per resulting row (row 2..6):
#   date (number), company, federal_tax_941, state_tax_edd, 941_payment_amount, EDD_payment_amount
$data = @(
    @(11012024, $companyName, "35783,00", "4438,24", "35783,00", "4438,24"),
    @(11082024, $companyName, "36233,65", "4416,39", "36233,65", "4416,39"),
    @(11152024, $companyName, "38771,93", "4793,73", "38771,93", "4793,73"),
    @(11222024, $companyName, "37915,74", "4702,90", "37915,74", "4702,90"),
    @(11292024, $companyName, "41179,84", "5151,17", "41179,84", "5151,26")
)

# Remove the now-superfluous rows (original sheet has 12 rows of data+header,
# the consolidated sheet only needs 6). Deleting from the bottom keeps the
# remaining row indices stable while we delete.
$ws.Rows("7:12").Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
